# Change field label from "Folder" to "Functional Area" in cell B1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Functional Area"

# Move the active selection to B1 (matches the saved sheet view state).
$ws.Range("B1").Select()
